# Weekly price update: insert a new week's Kiwi price entries (Primera/Segunda)
# at row 365, pushing the existing data down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 365-366 (existing rows 365+ shift down to 367+)
$ws.Rows("365:366").Insert()

# Row 365: new "Primera" entry for 2022-08-25
$ws.Cells.Item(365, 1).Value = 8
$ws.Cells.Item(365, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(365, 3).Value = "Coquimbo"
$ws.Cells.Item(365, 4).Value = (Get-Date -Year 2022 -Month 8 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(365, 5).Value = 4
$ws.Cells.Item(365, 6).Value = "Fruta"
$ws.Cells.Item(365, 7).Value = 100101
$ws.Cells.Item(365, 8).Value = "Berries"
$ws.Cells.Item(365, 9).Value = 100101007
$ws.Cells.Item(365, 10).Value = "Kiwi"
$ws.Cells.Item(365, 11).Value = "Hayward"
$ws.Cells.Item(365, 12).Value = "Primera"
$ws.Cells.Item(365, 13).Value = 16
$ws.Cells.Item(365, 14).Value = 210000
$ws.Cells.Item(365, 15).Value = 220000
$ws.Cells.Item(365, 16).Value = 215000
$ws.Cells.Item(365, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(365, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(365, 19).Value = 478
$ws.Cells.Item(365, 20).Value = 450

# Row 366: new "Segunda" entry for 2022-08-25
$ws.Cells.Item(366, 1).Value = 8
$ws.Cells.Item(366, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(366, 3).Value = "Coquimbo"
$ws.Cells.Item(366, 4).Value = (Get-Date -Year 2022 -Month 8 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(366, 5).Value = 4
$ws.Cells.Item(366, 6).Value = "Fruta"
$ws.Cells.Item(366, 7).Value = 100101
$ws.Cells.Item(366, 8).Value = "Berries"
$ws.Cells.Item(366, 9).Value = 100101007
$ws.Cells.Item(366, 10).Value = "Kiwi"
$ws.Cells.Item(366, 11).Value = "Hayward"
$ws.Cells.Item(366, 12).Value = "Segunda"
$ws.Cells.Item(366, 13).Value = 20
$ws.Cells.Item(366, 14).Value = 180000
$ws.Cells.Item(366, 15).Value = 190000
$ws.Cells.Item(366, 16).Value = 185000
$ws.Cells.Item(366, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(366, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(366, 19).Value = 411
$ws.Cells.Item(366, 20).Value = 450
